$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Copy($ws.Range("A6"))

$ws.Range("A6").Value = 42611.887002314812
$ws.Range("B6").Value = 26
$ws.Range("C6").Value = 62
$ws.Range("D6").Value = 35
$ws.Range("E6").Value = 66
$ws.Range("F6").Value = 33
$ws.Range("G6").Value = 21199
$ws.Range("H6").Value = 16025
$ws.Range("I6").Value = 881
$ws.Range("J6").Value = 190
$ws.Range("K6").Value = 108
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 3
$ws.Range("N6").Value = "Named"
